$d = $word.ActiveDocument

# Locate the whole "Mapky v tomto dokumente..." sentence (it spans three
# differently-formatted runs: plain text, the hyperlink-styled URL, and the
# trailing ").") and select it as a single Range.
$rng = $d.Content
$found = $rng.Find.Execute(
    "Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/).",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the astromap sentence to update"
}

$start = $rng.Start

# Remove the old runs entirely, then insert replacement markup at the same
# spot: an empty run followed by one plain run holding the updated sentence
# (year bumped from 2018 to 2022), matching how Word collapses a retyped
# selection into fresh, unformatted runs.
$rng.Delete()

$ins = $d.Range($start, $start)
$xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:t>Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).</w:t></w:r></w:p>"
$ins.InsertXML($xml)
